$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.888.58'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '3.428.92'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.18'
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.00'
$ws.Range("E6").Value = '  -1.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.607'
$ws.Range("E7").Value = '  +3.74%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '3.430.77'
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("E10").Value = '  -1.43%  '
$ws.Range("E11").Value = '  -1.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.441'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '4.022.44'
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("E15").Value = '  -3.32%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.75'
$ws.Range("E16").Value = '  -3.57%  '
$ws.Range("D17").Value = '64.886.02'
$ws.Range("D18").Value = '3.467.40'
$ws.Range("E18").Value = '  -2.01%  '
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("E20").Value = '  -2.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '381.23'
$ws.Range("E21").Value = '  -2.66%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.03'
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.550'
$ws.Range("E23").Value = '  -0.10%  '
$ws.Range("E24").Value = '  +0.19%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.40'
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("E26").Value = '  -4.44%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.12'
$ws.Range("E27").Value = '  +3.58%  '
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  +2.84%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.20'
$ws.Range("E31").Value = '  -3.79%  '
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.11'
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("E35").Value = '  +3.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.33'
$ws.Range("E36").Value = '  -1.84%  '
$ws.Range("E37").Value = '  -2.71%  '
$ws.Range("D38").Value = '2.935.13'
$ws.Range("E38").Value = '  -4.68%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0758'
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.75'
$ws.Range("E40").Value = '  +3.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.49'
$ws.Range("E41").Value = '  -2.61%  '
$ws.Range("E42").Value = '  +1.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.99'
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.02'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.773'
$ws.Range("E46").Value = '  -0.58%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '318.18'
$ws.Range("E47").Value = '  +1.90%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.26'
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.08'
$ws.Range("E49").Value = '  -4.07%  '
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("E51").Value = '  -2.13%  '
